# Locate the existing "all" bullet paragraph (3rd-level bullet list:
# in / any / all) and the blank paragraph that immediately follows it.
$d = $word.ActiveDocument

$allPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "all") {
        $allPara = $p
    }
}
$afterAll = $allPara.Next()

# Insert three new paragraphs right after that blank paragraph (i.e.
# before "Practical Demonstration"):
#   1) empty paragraph
#   2) paragraph with the new explanatory sentence (which also carries a
#      relocated "_GoBack" bookmark right after "operator the or")
#   3) empty paragraph
$afterAll.Range.InsertParagraphAfter()
$p1 = $afterAll.Next()

$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()

$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()

# --- paragraph 2: insert the sentence text ---
$sentence = "## All means 19 or 37 or 69 so in the any operator the or condition is taken"
$start = $p2.Range.Start
$ins = $d.Range($start, $start)
$ins.InsertAfter($sentence)

# Format the whole new sentence paragraph run(s): Segoe UI / #333333 / 10.5pt
$p2.Range.Font.Name = "Segoe UI"
$p2.Range.Font.NameFarEast = "Times New Roman"
$p2.Range.Font.Color = 3355443
$p2.Range.Font.Size = 10.5

# Move the document's "_GoBack" bookmark into the middle of the new
# sentence, right after "...operator the or" (Word keeps only one
# "_GoBack" bookmark at a time, so re-adding it here automatically
# removes it from its old location after "Products " in the 2nd SQL
# paragraph).
$markerOffset = $sentence.IndexOf("operator the or") + ("operator the or").Length
$bmPos = $start + $markerOffset
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Also apply the matching paragraph-mark formatting (Segoe UI / #333333 /
# 10.5pt) to the two blank paragraphs that were inserted alongside it.
foreach ($blank in @($p1, $p3)) {
    $blank.Range.Font.Name = "Segoe UI"
    $blank.Range.Font.NameFarEast = "Times New Roman"
    $blank.Range.Font.Color = 3355443
    $blank.Range.Font.Size = 10.5
}
